# Weekly update: insert two new price records (rows 143-144) for
# Feria Lagunitas de Puerto Montt - Ají, shifting the existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 143 (pushes current rows 143.. down to 145..)
$ws.Rows("143:144").Insert()

# New row 143
$ws.Cells.Item(143, 1).Value = 4
$ws.Cells.Item(143, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(143, 3).Value = "Los Lagos"
$ws.Cells.Item(143, 4).Value = 44505
$ws.Cells.Item(143, 5).Value = 10
$ws.Cells.Item(143, 6).Value = 100112021
$ws.Cells.Item(143, 7).Value = "Ají"
$ws.Cells.Item(143, 8).Value = "Inferno"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 70
$ws.Cells.Item(143, 11).Value = 33000
$ws.Cells.Item(143, 12).Value = 33000
$ws.Cells.Item(143, 13).Value = 33000
$ws.Cells.Item(143, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(143, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(143, 16).Value = 2750
$ws.Cells.Item(143, 17).Value = 12
$ws.Cells.Item(143, 18).Value = "Hortaliza"

# New row 144
$ws.Cells.Item(144, 1).Value = 4
$ws.Cells.Item(144, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(144, 3).Value = "Los Lagos"
$ws.Cells.Item(144, 4).Value = 44505
$ws.Cells.Item(144, 5).Value = 10
$ws.Cells.Item(144, 6).Value = 100112021
$ws.Cells.Item(144, 7).Value = "Ají"
$ws.Cells.Item(144, 8).Value = "Inferno"
$ws.Cells.Item(144, 9).Value = "Segunda"
$ws.Cells.Item(144, 10).Value = 70
$ws.Cells.Item(144, 11).Value = 28000
$ws.Cells.Item(144, 12).Value = 28000
$ws.Cells.Item(144, 13).Value = 28000
$ws.Cells.Item(144, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(144, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(144, 16).Value = 2333
$ws.Cells.Item(144, 17).Value = 12
$ws.Cells.Item(144, 18).Value = "Hortaliza"
